$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 80 (old rows 80..102 shift down to 82..104).
$ws.Rows.Item(80).Insert()
$ws.Rows.Item(80).Insert()

# --- New row 80 ---
$ws.Cells.Item(80, 1).Value = 9
$ws.Cells.Item(80, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(80, 3).Value = "Metropolitana"
$ws.Cells.Item(80, 4).Value = [DateTime]"2022-09-05"
$ws.Cells.Item(80, 5).Value = 13
$ws.Cells.Item(80, 6).Value = 100114002
$ws.Cells.Item(80, 7).Value = "Camote"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 900
$ws.Cells.Item(80, 11).Value = 14000
$ws.Cells.Item(80, 12).Value = 15000
$ws.Cells.Item(80, 13).Value = 14667
$ws.Cells.Item(80, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(80, 15).Value = "Perú"
$ws.Cells.Item(80, 16).Value = 815
$ws.Cells.Item(80, 17).Value = 18
$ws.Cells.Item(80, 18).Value = "Hortaliza"

# --- New row 81 ---
$ws.Cells.Item(81, 1).Value = 9
$ws.Cells.Item(81, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(81, 3).Value = "Metropolitana"
$ws.Cells.Item(81, 4).Value = [DateTime]"2022-09-05"
$ws.Cells.Item(81, 5).Value = 13
$ws.Cells.Item(81, 6).Value = 100114002
$ws.Cells.Item(81, 7).Value = "Camote"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 1300
$ws.Cells.Item(81, 11).Value = 12000
$ws.Cells.Item(81, 12).Value = 13000
$ws.Cells.Item(81, 13).Value = 12538
$ws.Cells.Item(81, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(81, 15).Value = "Perú"
$ws.Cells.Item(81, 16).Value = 697
$ws.Cells.Item(81, 17).Value = 18
$ws.Cells.Item(81, 18).Value = "Hortaliza"

Write-Output "done"
